$wb = $excel.ActiveWorkbook

# --- Overview sheet: handoff/handback status text for the dda7d53b row changed
#     from "Ready for handoff" to "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: same status text update, plus a new "Error Detail" message
#     in column P (row 3) describing the handback/handoff filename mismatch,
#     and widen column P (Error Detail) to fit the longer message.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: 5byrfwkf.poj is different with handoff file name: dda7d53b-03e4-4ba5-842d-32d3446d6980.f2608c516ae9e939592db6057fc0b4ada4181ac8.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: same status text update, plus a new "Error Detail" message
#     in column P (row 3), and widen column P to match.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: 5byrfwkf.poj is different with handoff file name: dda7d53b-03e4-4ba5-842d-32d3446d6980.f2608c516ae9e939592db6057fc0b4ada4181ac8.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
